$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.337.83"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.282.34"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "505.25"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "129.52"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.0956"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").Value = "0.334"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "2.687.89"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "22.95"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "54.291.12"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "2.285.33"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "10.29"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "4.13"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "305.25"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "61.83"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").Value = "7.34"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").Value = "174.57"
$ws.Range("E27").Value = "  +5.25%  "
$ws.Range("D28").Value = "1.63"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").Value = "6.03"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "0.0₃0691"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "1.09"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "17.79"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "0.966"
$ws.Range("E34").Value = "  +10.19%  "
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "3.76"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "4.86"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "125.18"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "0.0497"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "0.0894"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "0.548"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "241.73"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "0.0207"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "16.47"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("E51").Value = "  +0.20%  "
